# Changing isQuestion to messageType
#
# Column I used to hold a boolean ISQUESTION flag. It is replaced with a
# text MESSAGETYPE classification (question / message / final-message /
# splitting). Also fills in the previously-missing ID value in row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing ID for row 3 (IDs should run 1..9)
$ws.Range("A3").Value = 2

# Replace the ISQUESTION header with MESSAGETYPE
$ws.Range("I1").Value = "MESSAGETYPE"

# Replace the boolean ISQUESTION values with the new MESSAGETYPE values
$ws.Range("I2").Value = "question"
$ws.Range("I3").Value = "final-message"
$ws.Range("I4").Value = "message"
$ws.Range("I5").Value = "splitting"
$ws.Range("I6").Value = "message"
$ws.Range("I7").Value = "question"
$ws.Range("I8").Value = "final-message"
$ws.Range("I9").Value = "splitting"
$ws.Range("I10").Value = "message"

# Move the active selection to reflect where editing finished
[void]$ws.Range("I11").Select()
